$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 15287
$ws.Range("C6").Value = 27002685
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C12").Value = 87849554
$ws.Range("C13").Value = 6244281
